# Append a freshly-scraped Lancers listing at the top (row 3) of the
# "ランサーズ" sheet, push the previously-scraped rows down by one, refresh
# their scrape timestamp, and append one more brand-new listing as the new
# last row. Mirrors the site's "insert newest at top of this batch, shift
# the rest down, append overflow at bottom" scrape-merge behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-09-20 01:13:59"

# Final contents for every data row (row 2 is untouched by this commit).
$rows = @(
    @{ Row = 3;  A = $newTimestamp; B = "【急募】スマホアプリ自動化デモ開発(LLM連携)"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5397035"; G = 228; H = "★スマホアプリ ◆開発,自動化 ◇アプリ" },
    @{ Row = 4;  A = $newTimestamp; B = "システムの開発補助や運営サポート【フルリモート×長期】"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5323359"; G = 83; H = "◆開発" },
    @{ Row = 5;  A = $newTimestamp; B = "初回 【フィンテック/ブリッジ】金融資産管理システムの要件定義/仕様伝達/進捗管理(日/英|フルリモート)"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5396502"; G = 60; H = "◇管理" },
    @{ Row = 6;  A = $newTimestamp; B = "【フィンテック/QA】海外の金融資産管理システムのテスト設計・品質保証"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5396510"; G = 53; H = "◇管理" },
    @{ Row = 7;  A = $newTimestamp; B = "初回 【医療関連】会員制サイト構築のパートナーを探しています"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5396953"; G = 45; H = "◇サイト" },
    @{ Row = 8;  A = $newTimestamp; B = "会員情報サイトの新規構築プロジェクト"; C = "システム開発"; D = "1,000,000 円 ~ 3,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5396688"; G = 45; H = "◇サイト" },
    @{ Row = 9;  A = $newTimestamp; B = "【急募】PHP・Lalavelでの既存プログラム改修依頼"; C = "システム開発"; D = "100,000 円 ~ 200,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5396563"; G = 33; H = "○PHP" },
    @{ Row = 10; A = $newTimestamp; B = "【急募】WordPress記事をCoopelで自動投稿設定できる方を探しています!"; C = "システム開発"; D = "5,000 円 ~ 10,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5396609"; G = 25; H = "○WordPress" },
    @{ Row = 11; A = $newTimestamp; B = "【急募】データ統合基盤の設計・刷新プロジェクト"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5396957"; G = 25; H = $null },
    @{ Row = 12; A = $newTimestamp; B = "初回 Web広告のタグ設置・動作確認";                 C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定";      E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5397007"; G = 18; H = $null }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    if ($row.H -ne $null) {
        $ws.Range("H$r").Value = $row.H
    }
}

# Widen the "skill summary" column (H) to fit the new longer tag lists.
$ws.Columns.Item(8).ColumnWidth = 21.17

# Rebuild the URL hyperlinks for every data row (2..12) so rId1..rId11 line
# up positionally with F2..F12 again, each pointing at that row's own URL.
$ws.Range("F2:F12").Hyperlinks.Delete()
for ($r = 2; $r -le 12; $r++) {
    $target = $ws.Range("F$r").Value()
    $ws.Hyperlinks.Add($ws.Range("F$r"), $target)
    $ws.Range("F$r").Style = "Hyperlink"
}
